$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> (Ticket Sales [Q], Embarking [R]) new values
$updates = @{
    3   = @(13, 1)
    10  = @(39, 26)
    17  = @(1, 1)
    23  = @(14, 5)
    32  = @(49, 5)
    40  = @(61, 40)
    49  = @(45, 1)
    58  = @(8, 7)
    66  = @(11, 8)
    74  = @(95, 83)
    78  = @(78, 36)
    89  = @(65, 43)
    97  = @(58, 29)
    106 = @(59, 52)
    115 = @(74, 24)
    124 = @(22, 16)
    133 = @(81, 23)
    142 = @(11, 2)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("Q$row").Value = $vals[0]
    $ws.Range("R$row").Value = $vals[1]
}
